# Generate Report for Handback
#
# For the "aa987eea-c657-40b7-980a-57260e08d4da" file, a handback was
# processed but was found to be based on a stale version of the source
# document. Populate the per-language "Latest Target File" / "Latest
# Handback File" / "Latest Handback DateTime" / "Error Detail" columns
# (row 7) on both the zh-cn and de-de status sheets.

$wb = $excel.ActiveWorkbook

$latestUrl = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/43b628cc9be7f72f37391f7c8a40de4a57c36a2f/e2e/aa987eea-c657-40b7-980a-57260e08d4da.md"
$displayName = "aa987eea-c657-40b7-980a-57260e08d4da.md"
$errorDetail = "The version of handback file is not the latest, current: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/51a4e54d418a3f5ce762c00ee4e5bfa1d5bd0959/e2e/aa987eea-c657-40b7-980a-57260e08d4da.md, latest: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/43b628cc9be7f72f37391f7c8a40de4a57c36a2f/e2e/aa987eea-c657-40b7-980a-57260e08d4da.md."

# ---- zh-cn sheet ----
$wsZhCn = $wb.Worksheets.Item("zh-cn")

$iCell = $wsZhCn.Range("I7")
$wsZhCn.Hyperlinks.Add($iCell, $latestUrl, "", "", $displayName)

$wsZhCn.Range("J7").Value = "aa987eea-c657-40b7-980a-57260e08d4da.96c7cc6b202cd0b982d15a08563ae54cad495dda.zh-cn.xlf"
$wsZhCn.Range("K7").Value = "2016-08-27 20:53:52"
$wsZhCn.Range("P7").Value = $errorDetail

# ---- de-de sheet ----
$wsDeDe = $wb.Worksheets.Item("de-de")

$iCellDeDe = $wsDeDe.Range("I7")
$wsDeDe.Hyperlinks.Add($iCellDeDe, $latestUrl, "", "", $displayName)

$wsDeDe.Range("J7").Value = "aa987eea-c657-40b7-980a-57260e08d4da.96c7cc6b202cd0b982d15a08563ae54cad495dda.de-de.xlf"
$wsDeDe.Range("K7").Value = "2016-08-27 20:53:58"
$wsDeDe.Range("P7").Value = $errorDetail
